$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 31 (ALC)
$ws.Range("H31").Value = 866
$ws.Range("I31").Value = 866
$ws.Range("K31").Value = 2598
$ws.Range("M31").Value = -2368

# Row 129 (ALC)
$ws.Range("H129").Value = 25005214
$ws.Range("I129").Value = 83335630
$ws.Range("J129").Value = 6462
$ws.Range("K129").Value = 250006890
$ws.Range("L129").Value = 19386
$ws.Range("M129").Value = -250001890
$ws.Range("N129").Value = -29386

$ws = $wb.Worksheets.Item("ARM")
# Row 2 (ARM)
$ws.Range("H2").Value = 8930785
$ws.Range("I2").Value = 27778758
$ws.Range("J2").Value = 2798.2104
$ws.Range("K2").Value = 27778758
$ws.Range("L2").Value = 2798.2104
$ws.Range("M2").Value = -27778645
$ws.Range("N2").Value = -3024.2104

# Row 45 (ARM)
$ws.Range("H45").Value = 1796.3334
$ws.Range("J45").Value = 4732.3335
$ws.Range("L45").Value = 4732.3335
$ws.Range("N45").Value = -5486.3335

# Row 101 (ARM)
$ws.Range("H101").Value = 39701.5
$ws.Range("J101").Value = 39701.5
$ws.Range("L101").Value = 39701.5
$ws.Range("N101").Value = -46191.5

# Row 116 (ARM)
$ws.Range("H116").Value = 8930785
$ws.Range("I116").Value = 27778758
$ws.Range("J116").Value = 2798.2104
$ws.Range("K116").Value = 27778758
$ws.Range("L116").Value = 2798.2104
$ws.Range("M116").Value = -27776464
$ws.Range("N116").Value = -7386.2104

$ws = $wb.Worksheets.Item("BSM")
# Row 3 (BSM)
$ws.Range("H3").Value = 8930785
$ws.Range("I3").Value = 27778758
$ws.Range("J3").Value = 2798.2104
$ws.Range("K3").Value = 27778758
$ws.Range("L3").Value = 2798.2104
$ws.Range("M3").Value = -27778644
$ws.Range("N3").Value = -3026.2104

$ws = $wb.Worksheets.Item("CRP")
# Row 22 (CRP)
$ws.Range("H22").Value = 1923.875
$ws.Range("J22").Value = 4666.6665
$ws.Range("L22").Value = 4666.6665
$ws.Range("N22").Value = -5366.6665

# Row 31 (CRP)
$ws.Range("H31").Value = 1076997.8
$ws.Range("I31").Value = 1588242.5
$ws.Range("J31").Value = 3383.7666
$ws.Range("K31").Value = 1588242.5
$ws.Range("L31").Value = 3383.7666
$ws.Range("M31").Value = -1587947.5
$ws.Range("N31").Value = -3973.7666

# Row 34 (CRP)
$ws.Range("H34").Value = 1076997.8
$ws.Range("I34").Value = 1588242.5
$ws.Range("J34").Value = 3383.7666
$ws.Range("K34").Value = 1588242.5
$ws.Range("L34").Value = 3383.7666
$ws.Range("M34").Value = -1588040.5
$ws.Range("N34").Value = -3787.7666

# Row 58 (CRP)
$ws.Range("H58").Value = 17860258
$ws.Range("I58").Value = 2202.2104
$ws.Range("J58").Value = 55560600
$ws.Range("K58").Value = 2202.2104
$ws.Range("L58").Value = 55560600
$ws.Range("M58").Value = -1999.2104
$ws.Range("N58").Value = -55561006

# Row 99 (CRP)
$ws.Range("H99").Value = 2072.7334
$ws.Range("I99").Value = 1344.1
$ws.Range("J99").Value = 3530
$ws.Range("K99").Value = 1344.1
$ws.Range("L99").Value = 3530
$ws.Range("M99").Value = 153.9000000000001
$ws.Range("N99").Value = -6526

# Row 126 (CRP)
$ws.Range("H126").Value = 2072.7334
$ws.Range("I126").Value = 1344.1
$ws.Range("J126").Value = 3530
$ws.Range("K126").Value = 4032.3
$ws.Range("L126").Value = 10590
$ws.Range("M126").Value = -1562.3
$ws.Range("N126").Value = -15530

# Row 134 (CRP)
$ws.Range("H134").Value = 1837.4038
$ws.Range("I134").Value = 1457.95
$ws.Range("J134").Value = 3102.25
$ws.Range("K134").Value = 4373.85
$ws.Range("L134").Value = 9306.75
$ws.Range("M134").Value = -1838.85
$ws.Range("N134").Value = -14376.75

# Row 136 (CRP)
$ws.Range("H136").Value = 17860258
$ws.Range("I136").Value = 2202.2104
$ws.Range("J136").Value = 55560600
$ws.Range("K136").Value = 6606.6312
$ws.Range("L136").Value = 166681800
$ws.Range("M136").Value = -4056.6312
$ws.Range("N136").Value = -166686900

$ws = $wb.Worksheets.Item("CUL")
# Row 5 (CUL)
$ws.Range("H5").Value = 767.08826
$ws.Range("I5").Value = 452.86667
$ws.Range("J5").Value = 3123.75
$ws.Range("K5").Value = 1358.60001
$ws.Range("L5").Value = 9371.25
$ws.Range("M5").Value = -1246.60001
$ws.Range("N5").Value = -9595.25

# Row 55 (CUL)
$ws.Range("H55").Value = 3118.75
$ws.Range("J55").Value = 4000
$ws.Range("L55").Value = 12000
$ws.Range("N55").Value = -12354

# Row 122 (CUL)
$ws.Range("H122").Value = 741
$ws.Range("I122").Value = 371
$ws.Range("J122").Value = 1077.3636
$ws.Range("K122").Value = 3339
$ws.Range("L122").Value = 9696.2724
$ws.Range("M122").Value = -889
$ws.Range("N122").Value = -14596.2724

# Row 135 (CUL)
$ws.Range("H135").Value = 767.08826
$ws.Range("I135").Value = 452.86667
$ws.Range("J135").Value = 3123.75
$ws.Range("K135").Value = 4075.80003
$ws.Range("L135").Value = 28113.75
$ws.Range("M135").Value = -1540.80003
$ws.Range("N135").Value = -33183.75

$ws = $wb.Worksheets.Item("GSM")
# Row 119 (GSM)
$ws.Range("H119").Value = 20000
$ws.Range("J119").Value = 20000
$ws.Range("L119").Value = 20000
$ws.Range("N119").Value = -29676

# Row 126 (GSM)
$ws.Range("H126").Value = 2899.6
$ws.Range("I126").Value = 1558.3334
$ws.Range("J126").Value = 4137.6924
$ws.Range("K126").Value = 4675.0002
$ws.Range("L126").Value = 12413.0772
$ws.Range("M126").Value = -2205.0002
$ws.Range("N126").Value = -17353.0772

# Row 132 (GSM)
$ws.Range("H132").Value = 3770.5676
$ws.Range("I132").Value = 2853.1738
$ws.Range("J132").Value = 5277.7144
$ws.Range("K132").Value = 8559.5214
$ws.Range("L132").Value = 15833.1432
$ws.Range("M132").Value = -6029.5214
$ws.Range("N132").Value = -20893.1432

$ws = $wb.Worksheets.Item("LTW")
# Row 7 (LTW)
$ws.Range("H7").Value = 2467
$ws.Range("I7").Value = 1201.3334
$ws.Range("J7").Value = 3099.8333
$ws.Range("K7").Value = 1201.3334
$ws.Range("L7").Value = 3099.8333
$ws.Range("M7").Value = -1089.3334
$ws.Range("N7").Value = -3323.8333

# Row 40 (LTW)
$ws.Range("H40").Value = 3133.3333
$ws.Range("I40").Value = 800
$ws.Range("J40").Value = 3600
$ws.Range("K40").Value = 800
$ws.Range("L40").Value = 3600
$ws.Range("M40").Value = -664
$ws.Range("N40").Value = -3872

# Row 61 (LTW)
$ws.Range("H61").Value = 100004200
$ws.Range("I61").Value = 166669660
$ws.Range("J61").Value = 5999.75
$ws.Range("K61").Value = 166669660
$ws.Range("L61").Value = 5999.75
$ws.Range("M61").Value = -166669458
$ws.Range("N61").Value = -6403.75

# Row 113 (LTW)
$ws.Range("H113").Value = 100004200
$ws.Range("I113").Value = 166669660
$ws.Range("J113").Value = 5999.75
$ws.Range("K113").Value = 166669660
$ws.Range("L113").Value = 5999.75
$ws.Range("M113").Value = -166667490
$ws.Range("N113").Value = -10339.75

# Row 126 (LTW)
$ws.Range("H126").Value = 2467
$ws.Range("I126").Value = 1201.3334
$ws.Range("J126").Value = 3099.8333
$ws.Range("K126").Value = 3604.0002
$ws.Range("L126").Value = 9299.499899999999
$ws.Range("M126").Value = -1134.0002
$ws.Range("N126").Value = -14239.4999

$ws = $wb.Worksheets.Item("WVR")
# Row 132 (WVR)
$ws.Range("H132").Value = 275905.3
$ws.Range("I132").Value = 436956.22
$ws.Range("J132").Value = 11321.643
$ws.Range("K132").Value = 1310868.66
$ws.Range("L132").Value = 33964.929
$ws.Range("M132").Value = -1308338.66
$ws.Range("N132").Value = -39024.929
